# Update the fixed "Date Placeholder" field text from 8/20/2015 to
# 9/19/2015 everywhere it appears: the slide master, every slide layout
# (custom layout), and the notes master. Slides themselves are not
# affected (none of them carry a cached date/time field).

$p = $ppt.ActivePresentation
$oldDate = "8/20/2015"
$newDate = "9/19/2015"

function Update-DateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DateShape($p.SlideMaster)

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateShape($layouts.Item($j))
}

# Notes master
Update-DateShape($p.NotesMaster)
